$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 691.6896400000001
$ws.Range("I15").Value = 691.6896400000001
$ws.Range("K15").Value = 2075.06892
$ws.Range("M15").Value = -1906.06892
$ws.Range("H64").Value = 7450
$ws.Range("J64").Value = 7450
$ws.Range("L64").Value = 7450
$ws.Range("N64").Value = -7946
$ws.Range("H67").Value = 7450
$ws.Range("J67").Value = 7450
$ws.Range("L67").Value = 7450
$ws.Range("N67").Value = -9166
$ws.Range("H70").Value = 1844822.1
$ws.Range("J70").Value = 5993.4
$ws.Range("L70").Value = 17980.2
$ws.Range("N70").Value = -18520.2
$ws.Range("H73").Value = 1844822.1
$ws.Range("J73").Value = 5993.4
$ws.Range("L73").Value = 17980.2
$ws.Range("N73").Value = -19852.2
$ws.Range("H74").Value = 5750
$ws.Range("J74").Value = 5750
$ws.Range("L74").Value = 5750
$ws.Range("N74").Value = -7622
$ws.Range("H77").Value = 5750
$ws.Range("J77").Value = 5750
$ws.Range("L77").Value = 28750
$ws.Range("N77").Value = -38110
$ws.Range("H113").Value = 7222
$ws.Range("I113").Value = 4999.6
$ws.Range("K113").Value = 4999.6
$ws.Range("M113").Value = -1745.6
$ws.Range("H132").Value = 34487084
$ws.Range("J132").Value = 5743.25
$ws.Range("L132").Value = 17229.75
$ws.Range("N132").Value = -22289.75
$ws.Range("H137").Value = 2766.2307
$ws.Range("I137").Value = 2087
$ws.Range("K137").Value = 6261
$ws.Range("M137").Value = -3711
$ws.Range("H138").Value = 3266
$ws.Range("I138").Value = 2498.6
$ws.Range("J138").Value = 4225.25
$ws.Range("K138").Value = 7495.799999999999
$ws.Range("L138").Value = 12675.75
$ws.Range("M138").Value = -2355.799999999999
$ws.Range("N138").Value = -22955.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7338.9023
$ws.Range("I32").Value = 7497.4
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 7497.4
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -7210.4
$ws.Range("N32").Value = -1573
$ws.Range("H45").Value = 3733
$ws.Range("I45").Value = 3733
$ws.Range("K45").Value = 3733
$ws.Range("M45").Value = -3356
$ws.Range("H110").Value = 2008.1765
$ws.Range("I110").Value = 969.9167
$ws.Range("K110").Value = 969.9167
$ws.Range("M110").Value = 1075.0833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3998.6667
$ws.Range("I20").Value = 3998.6667
$ws.Range("K20").Value = 3998.6667
$ws.Range("M20").Value = -3751.6667
$ws.Range("H86").Value = 2659.5
$ws.Range("I86").Value = 2659.5
$ws.Range("K86").Value = 2659.5
$ws.Range("M86").Value = -1536.5
$ws.Range("H89").Value = 2659.5
$ws.Range("I89").Value = 2659.5
$ws.Range("K89").Value = 13297.5
$ws.Range("M89").Value = -7681.5
$ws.Range("H94").Value = 3300
$ws.Range("I94").Value = 3450
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 3450
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -2999
$ws.Range("N94").Value = -3902
$ws.Range("H107").Value = 1582.2
$ws.Range("J107").Value = 3222
$ws.Range("L107").Value = 3222
$ws.Range("N107").Value = -7062
$ws.Range("H134").Value = 5600.0356
$ws.Range("I134").Value = 5836.16
$ws.Range("J134").Value = 3632.3333
$ws.Range("K134").Value = 17508.48
$ws.Range("L134").Value = 10896.9999
$ws.Range("M134").Value = -14973.48
$ws.Range("N134").Value = -15966.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22643
$ws.Range("J28").Value = 22643
$ws.Range("L28").Value = 22643
$ws.Range("N28").Value = -23133
$ws.Range("H31").Value = 3908
$ws.Range("I31").Value = 2242
$ws.Range("K31").Value = 2242
$ws.Range("M31").Value = -1947
$ws.Range("H34").Value = 3908
$ws.Range("I34").Value = 2242
$ws.Range("K34").Value = 2242
$ws.Range("M34").Value = -2040
$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 10000
$ws.Range("K86").Value = 10000
$ws.Range("M86").Value = -8877
$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 10000
$ws.Range("K89").Value = 50000
$ws.Range("M89").Value = -44384
$ws.Range("H92").Value = 42085
$ws.Range("J92").Value = 42085
$ws.Range("L92").Value = 42085
$ws.Range("N92").Value = -47077
$ws.Range("H99").Value = 3333.3333
$ws.Range("I99").Value = 3333.3333
$ws.Range("K99").Value = 3333.3333
$ws.Range("M99").Value = -1835.3333
$ws.Range("H126").Value = 3333.3333
$ws.Range("I126").Value = 3333.3333
$ws.Range("K126").Value = 9999.999899999999
$ws.Range("M126").Value = -7529.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3271
$ws.Range("J129").Value = 3429.875
$ws.Range("L129").Value = 10289.625
$ws.Range("N129").Value = -20289.625
$ws.Range("H137").Value = 8155.3335
$ws.Range("J137").Value = 9386.4
$ws.Range("L137").Value = 28159.2
$ws.Range("N137").Value = -38359.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8577.777
$ws.Range("I70").Value = 7700
$ws.Range("J70").Value = 8687.5
$ws.Range("K70").Value = 7700
$ws.Range("L70").Value = 8687.5
$ws.Range("M70").Value = -7430
$ws.Range("N70").Value = -9227.5
$ws.Range("H73").Value = 8577.777
$ws.Range("I73").Value = 7700
$ws.Range("J73").Value = 8687.5
$ws.Range("K73").Value = 7700
$ws.Range("L73").Value = 8687.5
$ws.Range("M73").Value = -6764
$ws.Range("N73").Value = -10559.5
$ws.Range("H80").Value = 7401.625
$ws.Range("I80").Value = 6883
$ws.Range("J80").Value = 8068.4287
$ws.Range("K80").Value = 6883
$ws.Range("L80").Value = 8068.4287
$ws.Range("M80").Value = -5885
$ws.Range("N80").Value = -10064.4287
$ws.Range("H83").Value = 7401.625
$ws.Range("I83").Value = 6883
$ws.Range("J83").Value = 8068.4287
$ws.Range("K83").Value = 34415
$ws.Range("L83").Value = 40342.14350000001
$ws.Range("M83").Value = -29423
$ws.Range("N83").Value = -50326.14350000001
$ws.Range("H97").Value = 1542.7778
$ws.Range("I97").Value = 1485.625
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1485.625
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -989.625
$ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 2350.375
$ws.Range("I102").Value = 2500.4285
$ws.Range("K102").Value = 2500.4285
$ws.Range("M102").Value = -878.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 849.1667
$ws.Range("J22").Value = 999
$ws.Range("L22").Value = 999
$ws.Range("N22").Value = -1589
$ws.Range("H27").Value = 849.1667
$ws.Range("J27").Value = 999
$ws.Range("L27").Value = 999
$ws.Range("N27").Value = -1213
$ws.Range("H61").Value = 1697.1765
$ws.Range("I61").Value = 1514.8572
$ws.Range("K61").Value = 1514.8572
$ws.Range("M61").Value = -1312.8572
$ws.Range("H113").Value = 1697.1765
$ws.Range("I113").Value = 1514.8572
$ws.Range("K113").Value = 1514.8572
$ws.Range("M113").Value = 655.1428000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1295
$ws.Range("J62").Value = 1295
$ws.Range("L62").Value = 1295
$ws.Range("N62").Value = -2543
$ws.Range("H65").Value = 1295
$ws.Range("J65").Value = 1295
$ws.Range("L65").Value = 6475
$ws.Range("N65").Value = -12715
$ws.Range("H107").Value = 374
$ws.Range("I107").Value = 365.66666
$ws.Range("K107").Value = 1096.99998
$ws.Range("M107").Value = 823.0000199999999
$ws.Range("H132").Value = 849
$ws.Range("I132").Value = 849
$ws.Range("K132").Value = 2547
$ws.Range("M132").Value = -17
